$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Directorio")

# The hierarchical ID / ID_Padre numbering scheme is being made consistent
# by appending a trailing "." to every node id (e.g. "1.1" -> "1.1.",
# "1" -> "1.", "2" -> "2."). Some of the existing values (e.g. "1.1.3.")
# already end in "." and are untouched.
#
# A handful of the new values (the ones that still look like a plain
# number with one trailing dot, e.g. "1.", "2.", "3.") would normally be
# re-interpreted by Excel's automatic number detection when assigned via
# .Value, which would also mint a brand-new cell style. To avoid that, we
# round-trip those specific values through a text-formula + paste-values
# special, which keeps the literal text (stored as a shared string) while
# preserving the cell's existing style.

function Set-TextLiteral($range, [string]$text) {
    $helper = $ws.Range("Z100")
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
    $helper.ClearContents()
}

# Column A = ID, Column B = ID_Padre
$ws.Range("A3").Value = "1.1."
Set-TextLiteral $ws.Range("B3") "1."

$ws.Range("A4").Value = "1.1.1."
$ws.Range("B4").Value = "1.1."

$ws.Range("A6").Value = "1.1.2."
$ws.Range("B6").Value = "1.1."

$ws.Range("A7").Value = "1.1.3."
$ws.Range("B7").Value = "1.1."

$ws.Range("A8").Value = "1.2."
Set-TextLiteral $ws.Range("B8") "1."

$ws.Range("A9").Value = "1.3."
Set-TextLiteral $ws.Range("B9") "1."

Set-TextLiteral $ws.Range("A10") "2."

Set-TextLiteral $ws.Range("A11") "3."

# Update the active cell selection (was D5, now D6)
$ws.Range("D6").Select()
